$wb = $excel.ActiveWorkbook

# --- "survey" sheet -------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# The old 3-row "inline yes_no" demo (i1/i2/i3, all labelled "Choose one:")
# is being turned into a "Which devices do you use?" note followed by a
# 4-row "inline yes_no" demo (desktop/laptop/smartphone/tablet). That is a
# net +2 rows, so insert 2 blank rows before the old "end screen" row (13).
$survey.Rows.Item(13).Insert()
$survey.Rows.Item(13).Insert()

# Row 10: was "select_one yes_no" / inline / i1 / "Choose one:" -> becomes
# a plain note introducing the new question.
$survey.Range("B10").Value = "note"
$survey.Range("C10").ClearContents()
$survey.Range("F10").ClearContents()
$survey.Range("G10").Value = "Which devices do you use?"

# Row 11: i2 -> desktop
$survey.Range("F11").Value = "desktop"
$survey.Range("G11").Value = "Desktop computer"

# Row 12: i3 -> laptop
$survey.Range("F12").Value = "laptop"
$survey.Range("G12").Value = "Laptop computer"

# Row 13 (new): smartphone
$survey.Range("B13").Value = "select_one yes_no"
$survey.Range("C13").Value = "inline"
$survey.Range("F13").Value = "smartphone"
$survey.Range("G13").Value = "Smartphone"

# Row 14 (new): tablet
$survey.Range("B14").Value = "select_one yes_no"
$survey.Range("C14").Value = "inline"
$survey.Range("F14").Value = "tablet"
$survey.Range("G14").Value = "Tablet"

# The content-provider-query row (now row 21, was row 19) gains an explicit
# "condition" value of FALSE.
$survey.Range("E21").Value = $false

# Column E ("condition") is narrower in the refactored form.
$survey.Columns.Item(5).ColumnWidth = 33.29

# --- "queries" sheet --------------------------------------------------------
$queries = $wb.Worksheets.Item("queries")
$queries.Range("C6").Value = "context"
